$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 21:50"
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 238823
$ws.Range("C4").Value = 23820
$ws.Range("D4").Value = 10360
$ws.Range("E4").Value = 222695
$ws.Range("F4").Value = 5421
$ws.Range("G4").Value = 666
$ws.Range("H4").Value = 5768

$ws.Range("A7").Value = "Alemania"
$ws.Range("B7").Value = 84636
$ws.Range("C7").Value = 6655
$ws.Range("D7").Value = 21400
$ws.Range("E7").Value = 62137
$ws.Range("F7").Value = 3936
$ws.Range("G7").Value = 168
$ws.Range("H7").Value = 1099

$ws.Range("A21").Value = "Israel"
$ws.Range("B21").Value = 6857
$ws.Range("C21").Value = 765
$ws.Range("D21").Value = 338
$ws.Range("E21").Value = 6483
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 10
$ws.Range("H21").Value = 36

$ws.Range("A23").Value = "Noruega"
$ws.Range("B23").Value = 5140
$ws.Range("C23").Value = 263
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = 5058
$ws.Range("F23").Value = 96
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 50

$ws.Range("A24").Value = "Australia"
$ws.Range("B24").Value = 5139
$ws.Range("C24").Value = 91
$ws.Range("D24").Value = 345
$ws.Range("E24").Value = 4767
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 27

$ws.Range("A32").Value = "Polonia"
$ws.Range("B32").Value = 2946
$ws.Range("C32").Value = 392
$ws.Range("D32").Value = 56
$ws.Range("E32").Value = 2833
$ws.Range("F32").Value = 50
$ws.Range("G32").Value = 14
$ws.Range("H32").Value = 57

$ws.Range("A33").Value = "Rumania"
$ws.Range("B33").Value = 2738
$ws.Range("C33").Value = 278
$ws.Range("D33").Value = 267
$ws.Range("E33").Value = 2356
$ws.Range("F33").Value = 78
$ws.Range("G33").Value = 23
$ws.Range("H33").Value = 115

$ws.Range("A37").Value = "Pakistan"
$ws.Range("B37").Value = 2386
$ws.Range("C37").Value = 268
$ws.Range("D37").Value = 107
$ws.Range("E37").Value = 2245
$ws.Range("F37").Value = 9
$ws.Range("G37").Value = 7
$ws.Range("H37").Value = 34

$ws.Range("A45").Value = "Peru"
$ws.Range("B45").Value = 1414
$ws.Range("C45").Value = 91
$ws.Range("D45").Value = 537
$ws.Range("E45").Value = 822
$ws.Range("F45").Value = 49
$ws.Range("G45").Value = 17
$ws.Range("H45").Value = 55

$ws.Range("A103").Value = "Costa de Marfil"
$ws.Range("B103").Value = 194
$ws.Range("C103").Value = 4
$ws.Range("D103").Value = 9
$ws.Range("E103").Value = 184
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 1

$ws.Range("A158").Value = "Nueva Caledonia"
$ws.Range("B158").Value = 18
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 1
$ws.Range("E158").Value = 17
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0

$ws.Range("A159").Value = "Gabon"
$ws.Range("B159").Value = 18
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 0
$ws.Range("E159").Value = 17
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 1

$ws.Range("A165").Value = "Namibia"
$ws.Range("B165").Value = 14
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 2
$ws.Range("E165").Value = 12
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0

$ws.Range("A166").Value = "Mongolia"
$ws.Range("B166").Value = 14
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 2
$ws.Range("E166").Value = 12
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0

$ws.Range("A167").Value = "Santa Lucia"
$ws.Range("B167").Value = 13
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 1
$ws.Range("E167").Value = 12
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

$ws.Range("A168").Value = "Benin"
$ws.Range("B168").Value = 13
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 1
$ws.Range("E168").Value = 12
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

$ws.Range("A171").Value = "Libia"
$ws.Range("B171").Value = 10
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 0
$ws.Range("E171").Value = 10
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

$ws.Range("A173").Value = "Mozambique"
$ws.Range("B173").Value = 10
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 10
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("A174").Value = "Seychelles"
$ws.Range("B174").Value = 10
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

$ws.Range("A175").Value = "Granada"
$ws.Range("B175").Value = 10
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 10
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

$ws.Range("A176").Value = "Laos"
$ws.Range("B176").Value = 10
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

$ws.Range("A179").Value = "Guinea-Bisau"
$ws.Range("B179").Value = 9
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 9
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

$ws.Range("A180").Value = "Suazilandia"
$ws.Range("B180").Value = 9
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 9
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0

$ws.Range("A181").Value = "Zimbabue"
$ws.Range("B181").Value = 9
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 8
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 1

$ws.Range("A182").Value = "Montserrat"
$ws.Range("B182").Value = 9
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 7
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 2

$ws.Range("A183").Value = "San Cristobal y Nieves"
$ws.Range("B183").Value = 8
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 8
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

